$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01134666666666667
$ws.Range("H2").Value = 0.03404
$ws.Range("I2").Value = 0.001209510404472147
$ws.Range("J2").Value = 0.001209510404472147
$ws.Range("M2").Value = 0.2799683333333333
$ws.Range("N2").Value = 0.839905
$ws.Range("O2").Value = 0.0294305463214559
$ws.Range("P2").Value = 0.0294305463214559
$ws.Range("Q2").Value = 0.003176707355555555
$ws.Range("R2").Value = 0.0285903662
$ws.Range("S2").Value = 0.00003559655198510038
$ws.Range("T2").Value = 0.00003559655198510038

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01134666666666667
$ws.Range("H3").Value = 0.03404
$ws.Range("I3").Value = 0.001209510404472147
$ws.Range("J3").Value = 0.001209510404472147
$ws.Range("O3").Value = 0.2486942046732164
$ws.Range("P3").Value = 0.2486942046732163
$ws.Range("Q3").Value = 0.02684383431555555
$ws.Range("R3").Value = 0.24159450884
$ws.Range("S3").Value = 0.0003007982280841808
$ws.Range("T3").Value = 0.0003007982280841807

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01134666666666667
$ws.Range("H4").Value = 0.03404
$ws.Range("I4").Value = 0.001209510404472147
$ws.Range("J4").Value = 0.001209510404472147
$ws.Range("M4").Value = 6.86709
$ws.Range("N4").Value = 20.60127
$ws.Range("O4").Value = 0.7218752490053277
$ws.Range("P4").Value = 0.7218752490053277
$ws.Range("Q4").Value = 0.0779185812
$ws.Range("R4").Value = 0.7012672308
$ws.Range("S4").Value = 0.0008731156244028656
$ws.Range("T4").Value = 0.0008731156244028656

# Row 5
$ws.Range("I5").Value = 0.8865539289740954
$ws.Range("J5").Value = 0.8865539289740952
$ws.Range("M5").Value = 0.2799683333333333
$ws.Range("N5").Value = 0.839905
$ws.Range("O5").Value = 0.0294305463214559
$ws.Range("P5").Value = 0.0294305463214559
$ws.Range("Q5").Value = 2.328481323397778
$ws.Range("R5").Value = 20.95633191058
$ws.Range("S5").Value = 0.02609176647314083
$ws.Range("T5").Value = 0.02609176647314083

# Row 6
$ws.Range("I6").Value = 0.8865539289740954
$ws.Range("J6").Value = 0.8865539289740952
$ws.Range("O6").Value = 0.2486942046732164
$ws.Range("P6").Value = 0.2486942046732163
$ws.Range("S6").Value = 0.2204808242661278
$ws.Range("T6").Value = 0.2204808242661278

# Row 7
$ws.Range("I7").Value = 0.8865539289740954
$ws.Range("J7").Value = 0.8865539289740952
$ws.Range("M7").Value = 6.86709
$ws.Range("N7").Value = 20.60127
$ws.Range("O7").Value = 0.7218752490053277
$ws.Range("P7").Value = 0.7218752490053277
$ws.Range("Q7").Value = 57.11321212908001
$ws.Range("R7").Value = 514.01890916172
$ws.Range("S7").Value = 0.6399813382348267
$ws.Range("T7").Value = 0.6399813382348266

# Row 8
$ws.Range("G8").Value = 1.052914333333334
$ws.Range("H8").Value = 3.158743
$ws.Range("I8").Value = 0.1122365606214325
$ws.Range("J8").Value = 0.1122365606214325
$ws.Range("M8").Value = 0.2799683333333333
$ws.Range("N8").Value = 0.839905
$ws.Range("O8").Value = 0.0294305463214559
$ws.Range("P8").Value = 0.0294305463214559
$ws.Range("Q8").Value = 0.2947826710461112
$ws.Range("R8").Value = 2.653044039415
$ws.Range("S8").Value = 0.003303183296329963
$ws.Range("T8").Value = 0.003303183296329963

# Row 9
$ws.Range("G9").Value = 1.052914333333334
$ws.Range("H9").Value = 3.158743
$ws.Range("I9").Value = 0.1122365606214325
$ws.Range("J9").Value = 0.1122365606214325
$ws.Range("O9").Value = 0.2486942046732164
$ws.Range("P9").Value = 0.2486942046732163
$ws.Range("Q9").Value = 2.490974551628112
$ws.Range("R9").Value = 22.418770964653
$ws.Range("S9").Value = 0.0279125821790044
$ws.Range("T9").Value = 0.02791258217900439

# Row 10
$ws.Range("G10").Value = 1.052914333333334
$ws.Range("H10").Value = 3.158743
$ws.Range("I10").Value = 0.1122365606214325
$ws.Range("J10").Value = 0.1122365606214325
$ws.Range("M10").Value = 6.86709
$ws.Range("N10").Value = 20.60127
$ws.Range("O10").Value = 0.7218752490053277
$ws.Range("P10").Value = 0.7218752490053277
$ws.Range("Q10").Value = 7.230457489290002
$ws.Range("R10").Value = 65.07411740361
$ws.Range("S10").Value = 0.08102079514609815
$ws.Range("T10").Value = 0.08102079514609815
